$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "67.127.04"
Set-TextCell $ws.Range("E2") "  +1.91%  "

# Row 3
Set-TextCell $ws.Range("D3") "3.907.73"
Set-TextCell $ws.Range("E3") "  +3.37%  "

# Row 4
Set-TextCell $ws.Range("D4") "0.999"
Set-TextCell $ws.Range("E4") "  +0.13%  "

# Row 5
Set-TextCell $ws.Range("D5") "465.00"
Set-TextCell $ws.Range("E5") "  +8.73%  "

# Row 6
Set-TextCell $ws.Range("D6") "144.27"
Set-TextCell $ws.Range("E6") "  +4.50%  "

# Row 7
Set-TextCell $ws.Range("D7") "0.629"
Set-TextCell $ws.Range("E7") "  +0.40%  "

# Row 8
Set-TextCell $ws.Range("D8") "0.998"
Set-TextCell $ws.Range("E8") "  -0.07%  "

# Row 9
Set-TextCell $ws.Range("D9") "0.741"
Set-TextCell $ws.Range("E9") "  +0.31%  "

# Row 10
Set-TextCell $ws.Range("D10") "0.166"
Set-TextCell $ws.Range("E10") "  +8.21%  "

# Row 11
Set-TextCell $ws.Range("D11") "0.0000342"
Set-TextCell $ws.Range("E11") "  +8.37%  "

# Row 12
Set-TextCell $ws.Range("D12") "43.04"
Set-TextCell $ws.Range("E12") "  +0.40%  "

# Row 13
Set-TextCell $ws.Range("D13") "10.43"
Set-TextCell $ws.Range("E13") "  -0.83%  "

# Row 14
Set-TextCell $ws.Range("D14") "4.519.57"
Set-TextCell $ws.Range("E14") "  +3.41%  "

# Row 15
Set-TextCell $ws.Range("D15") "15.28"
Set-TextCell $ws.Range("E15") "  +2.18%  "

# Row 16
Set-TextCell $ws.Range("D16") "3.916.60"
Set-TextCell $ws.Range("E16") "  +3.97%  "

# Row 17
Set-TextCell $ws.Range("E17") "  -0.37%  "

# Row 18
Set-TextCell $ws.Range("D18") "20.00"
Set-TextCell $ws.Range("E18") "  +0.12%  "

# Row 19
Set-TextCell $ws.Range("D19") "1.16"
Set-TextCell $ws.Range("E19") "  +2.31%  "

# Row 20
Set-TextCell $ws.Range("D20") "67.322.29"
Set-TextCell $ws.Range("E20") "  +2.13%  "

# Row 21
Set-TextCell $ws.Range("D21") "434.12"
Set-TextCell $ws.Range("E21") "  +6.97%  "

# Row 22
Set-TextCell $ws.Range("D22") "14.79"
Set-TextCell $ws.Range("E22") "  -2.91%  "

# Row 23
Set-TextCell $ws.Range("D23") "3.37"
Set-TextCell $ws.Range("E23") "  +3.64%  "

# Row 24
Set-TextCell $ws.Range("D24") "88.80"
Set-TextCell $ws.Range("E24") "  +4.35%  "

# Row 25
Set-TextCell $ws.Range("D25") "38.78"
Set-TextCell $ws.Range("E25") "  +5.79%  "

# Row 26
Set-TextCell $ws.Range("D26") "3.53"
Set-TextCell $ws.Range("E26") "  +7.68%  "

# Row 27
Set-TextCell $ws.Range("D27") "5.77"
Set-TextCell $ws.Range("E27") "  +6.70%  "

# Row 28
Set-TextCell $ws.Range("D28") "10.12"
Set-TextCell $ws.Range("E28") "  +1.91%  "

# Row 29
Set-TextCell $ws.Range("E29") "  -2.55%  "

# Row 30
Set-TextCell $ws.Range("D30") "736.83"
Set-TextCell $ws.Range("E30") "  +5.06%  "

# Row 31
Set-TextCell $ws.Range("D31") "13.71"
Set-TextCell $ws.Range("E31") "  -1.82%  "

# Row 32
Set-TextCell $ws.Range("E32") "  +0.99%  "

# Row 33
Set-TextCell $ws.Range("D33") "2.80"
Set-TextCell $ws.Range("E33") "  +0.70%  "

# Row 34
Set-TextCell $ws.Range("D34") "43.23"
Set-TextCell $ws.Range("E34") "  +6.45%  "

# Row 35
Set-TextCell $ws.Range("E35") "  +5.08%  "

# Row 36
Set-TextCell $ws.Range("D36") "58.17"
Set-TextCell $ws.Range("E36") "  +2.92%  "

# Row 37
Set-TextCell $ws.Range("E37") "  -0.04%  "

# Row 38
Set-TextCell $ws.Range("D38") "0.0₃0806"
Set-TextCell $ws.Range("E38") "  +18.97%  "

# Row 39
Set-TextCell $ws.Range("D39") "5.39"
Set-TextCell $ws.Range("E39") "  -6.67%  "

# Row 40
Set-TextCell $ws.Range("E40") "  +14.17%  "

# Row 41
Set-TextCell $ws.Range("D41") "0.0477"
Set-TextCell $ws.Range("E41") "  +0.76%  "

# Row 42
Set-TextCell $ws.Range("D42") "0.141"
Set-TextCell $ws.Range("E42") "  -1.52%  "

# Row 43
Set-TextCell $ws.Range("E43") "  +0.24%  "

# Row 44
Set-TextCell $ws.Range("D44") "0.335"
Set-TextCell $ws.Range("E44") "  +5.05%  "

# Row 45
Set-TextCell $ws.Range("E45") "  +5.08%  "

# Row 46
Set-TextCell $ws.Range("D46") "2.18"
Set-TextCell $ws.Range("E46") "  +5.42%  "

# Row 47
Set-TextCell $ws.Range("D47") "3.41"
Set-TextCell $ws.Range("E47") "  +1.88%  "

# Row 48
Set-TextCell $ws.Range("D48") "2.51"
Set-TextCell $ws.Range("E48") "  -2.92%  "

# Row 49 (was Stacks, now ApeXProtocol)
Set-TextCell $ws.Range("B49") "ApeXProtocol"
Set-TextCell $ws.Range("C49") "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextCell $ws.Range("D49") "3.13"
Set-TextCell $ws.Range("E49") "  -0.35%  "

# Row 50 (was ApeXProtocol, now Stacks)
Set-TextCell $ws.Range("B50") "Stacks"
Set-TextCell $ws.Range("C50") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell $ws.Range("D50") "2.91"
Set-TextCell $ws.Range("E50") "  +3.62%  "

# Row 51 (Monero price/volume update)
Set-TextCell $ws.Range("D51") "143.52"
Set-TextCell $ws.Range("E51") "  +1.31%  "
